$d = $word.ActiveDocument

function Append-Run($para, [string]$text, [bool]$italic) {
    $insertAt = $para.Range.End - 1
    $r = $d.Range($insertAt, $insertAt)
    $r.InsertAfter($text)
    if ($italic) {
        $r.Font.Italic = $true
    }
}

# Update the letter date in the first-page header (header3.xml, rId12, w:type="first")
$sec = $d.Sections(1)
$hdr = $sec.Headers.Item(2)
$hdr.Range.Find.Execute("2023-09-13", $true, $false, $false, $false, $false, $true, 1, $false, "2023-09-15", 2)

# Anchor: last paragraph currently in the document body (the "BILAGA 1 - Fridlysta arter" title)
$anchor = $d.Paragraphs.Last.Range

# --- new paragraph 0 (style=Heading1) ---
$anchor.InsertParagraphAfter()
$p0 = $d.Paragraphs.Last
$p0.Style = 'Heading1'
$anchor = $p0.Range

# --- new paragraph 1 (style=None) ---
$anchor.InsertParagraphAfter()
$p1 = $d.Paragraphs.Last
$p1.Style = 'Normal'
$anchor = $p1.Range

# --- new paragraph 2 (style=None) ---
$anchor.InsertParagraphAfter()
$p2 = $d.Paragraphs.Last
$p2.Style = 'Normal'
$anchor = $p2.Range

# --- new paragraph 3 (style=None) ---
$anchor.InsertParagraphAfter()
$p3 = $d.Paragraphs.Last
$p3.Style = 'Normal'
$anchor = $p3.Range

# --- new paragraph 4 (style=None) ---
$anchor.InsertParagraphAfter()
$p4 = $d.Paragraphs.Last
$p4.Style = 'Normal'
$anchor = $p4.Range

# --- new paragraph 5 (style=None) ---
$anchor.InsertParagraphAfter()
$p5 = $d.Paragraphs.Last
$p5.Style = 'Normal'
$anchor = $p5.Range

# --- new paragraph 6 (style=Heading2) ---
$anchor.InsertParagraphAfter()
$p6 = $d.Paragraphs.Last
$p6.Style = 'Heading2'
$anchor = $p6.Range

# --- new paragraph 7 (style=None) ---
$anchor.InsertParagraphAfter()
$p7 = $d.Paragraphs.Last
$p7.Style = 'Normal'
$anchor = $p7.Range

# --- new paragraph 8 (style=None) ---
$anchor.InsertParagraphAfter()
$p8 = $d.Paragraphs.Last
$p8.Style = 'Normal'
$anchor = $p8.Range

# --- new paragraph 9 (style=None) ---
$anchor.InsertParagraphAfter()
$p9 = $d.Paragraphs.Last
$p9.Style = 'Normal'
$anchor = $p9.Range

# --- new paragraph 10 (style=None) ---
$anchor.InsertParagraphAfter()
$p10 = $d.Paragraphs.Last
$p10.Style = 'Normal'
$anchor = $p10.Range

# --- new paragraph 11 (style=None) ---
$anchor.InsertParagraphAfter()
$p11 = $d.Paragraphs.Last
$p11.Style = 'Normal'
$anchor = $p11.Range

# --- new paragraph 12 (style=None) ---
$anchor.InsertParagraphAfter()
$p12 = $d.Paragraphs.Last
$p12.Style = 'Normal'
$anchor = $p12.Range

# --- fill text for paragraph 0 ---
Append-Run $p0 'Knärot – ekologi samt krav på livsmiljön' $false

# --- fill text for paragraph 1 ---
Append-Run $p1 'Knärot är fridlyst enligt 8 och 15 §§ artskyddsförordningen och klassad som sårbar (VU) enligt rödlistan 2020. Knärot är beroende av hög och jämn luftfuktighet i gamla, ostörda skogsmiljöer och är känslig för snabba förändringar av ljus-/vindförhållanden eller uttorkning. På grund av ett alltför intensivt skogsbruk har den minskat med 40 (25-50) % under de senaste 60 åren och i framtiden bedöms minskningstakten uppgå till 30 (20-40) %. Till följd av att arten har en dokumenterat högre minskningstakt iförhållande till sin generationstid än vad som tidigare varit känt (data från Riksskogstaxeringen) höjdes den till hotkategori sårbar (VU) i rödlistan 2020 (Artdatabanken, 2021).' $false

# --- fill text for paragraph 2 ---
Append-Run $p2 'Samuel Johnsons doktorsavhandling ' $false
Append-Run $p2 '“Retention Forestry as a Conservation Measure for Boreal Forest Ground Vegetation“' $true
Append-Run $p2 ' (SLU, Uppsala 2014) visar att det krävs väl tilltagna skyddszoner för att knärotens växtplatser inte ska ta skada av skogsbruksåtgärder i intilliggande områden: ' $false
Append-Run $p2 '“Study III shows that retention patches smaller than 0.5 ha do not lifeboat the sensitive forest herb G. repens, a species that depend on stable microclimatic conditions typical for intact forest stands.” ' $true
Append-Run $p2 'Vidare ' $false
Append-Run $p2 '“More sensitive forest species are not lifeboated in retention patches ranging from 0.05 to 0.5 ha (Papers II & III).”' $true

# --- fill text for paragraph 3 ---
Append-Run $p3 'Johnsons (2014) rekommendation på minst 50 meters breda skyddszoner runt knärotens växtplatser motsvarar en areal på 0,78 hektar, vilket ligger i linje med andra studier som gjorts på känsliga skogsarter: ' $false
Append-Run $p3 '“In study III I also show that translocated specimens of G. repens survives well in mature forests at least 50 m from the nearest edge to an open area. Moreover, measures of temperature and humidity show that such distances from an open area is far enough to offer a microclimate that is more stable compared to what present in retention patches of around 0.1 ha. This means that the very centre of a circular patch with radius 50 m (equals a size of 0.78 ha) should offer conditions similar to interior forest and would perhaps be a suitable habitat for G. repens and similar species. Previous studies from both North America and Sweden have also concluded that patches between 0.5 and one ha are sufficient for preserving interior forest vegetation as well as sensitive lichens and bryophytes (de Graaf & Roberts 2009; Halpern et al. 2012; Rudolphi et al. 2014).”' $true

# --- fill text for paragraph 4 ---
Append-Run $p4 'En nyligen publicerad vetenskaplig uppsats av Koelmeijer m.fl. (2022) inkluderar orkidén knärots skyddsbehov. I uppsatsen berörs problemet med uttorkning för växter, bl.a. för knärot, ett problem som blivit accentuerat på grund av den pågående klimatförändringen och torra somrar, t.ex. den exceptionellt torra sommaren 2018. I uppsatsen undersöks områden med tre olika avstånd från kalhyggeskant med avseende på skydd bl.a. för knärot. Det första området har avstånd upp till 20 m från hyggeskant (Strong edge effect), det andra 20 – 40 m från hyggeskant (Weak edge effect) och det tredje avser större avstånd från hyggeskant, där kanteffekten anses vara försumbar (Interior). Ett resultat var att man fann stor eller mycket stor uttorkningseffekt på känsliga och rödlistade skogsarter vid de kortare avstånden till hyggeskant, medan effekt av uttorkning inte konstaterades på större avstånd (Interior). För orkidén knärot fann man en rik förekomst (upp till 0,06 dm2/m2) på stort avstånd från hyggeskant (Interior), medan förekomsten var liten eller närmast försumbar i de områden som klassificerades som Weak edge effect respektive Strong edge effect. Arbetet påpekar att de allt oftare förekommande torra somrarna ger ytterligare skäl att utöka skyddsavståndet från hyggen till den fuktkrävande arten knärot (Koelmeijer m.fl., 2022).' $false

# --- fill text for paragraph 5 ---
Append-Run $p5 'Även Skogsstyrelsens egen vägledning för hänsyn till knärot ligger i linje med ovanstående forskningsstudier. Av vägledningen framgår det att för med hög sannolikhet kunna bevara befintliga förekomster krävs relativt stora avsättningar av uppvuxen skog med slutet och relativt tätt kronskikt. Som riktlinje kan krävas ett avstånd på 50 meter in från brynet för att vidmakthålla ett fungerande mikroklimat. Detta innebär att fristående hänsynsytor för många arter (kärlväxter, lavar och mossor) kan behöva ha en area överstigande 0,8 hektar (cirkelyta med radien 50 meter = 0,78 hektar) för att bibehålla lokalklimatet. Även ganska små förändringar i form av förändrade ljus- och fuktighetsförhållanden, till exempel till följd av gallring, kan leda till att arten försvinner till följd av konkurrens med mera ljuskrävande och snabbväxande arter (Skogsstyrelsen, 2022).' $false

# --- fill text for paragraph 6 ---
Append-Run $p6 'Referenser - knärot' $false

# --- fill text for paragraph 7 ---
Append-Run $p7 'de Graaf M & Roberts M.R., 2009. ' $false
Append-Run $p7 'Short-term response of the herbaceous layer within leave patches after harvest. ' $true
Append-Run $p7 'Forest Ecology and Management 257, 1014-1025' $false

# --- fill text for paragraph 8 ---
Append-Run $p8 'Halpern, C. B., Halaj, J., Evans, S. A., & Dovciak, M., 2012. ' $false
Append-Run $p8 'Level and pattern of overstory retention interact to shape long-term responses of understories to timber harvest. ' $true
Append-Run $p8 'Ecological Applications, 22, 2049-2064 ' $false

# --- fill text for paragraph 9 ---
Append-Run $p9 'Koelmeijer, I. A., Ehrlén, J., Jönsson, M., De Frenne, P., Berg, P., Andersson, J., Weibull, H. & Hylander, N. 2022. ' $false
Append-Run $p9 'Interactive effects of drought and edge exposure on old-growth forest understory species. ' $true
Append-Run $p9 'Landscape Ecology, 37, sid 1839-1853' $false

# --- fill text for paragraph 10 ---
Append-Run $p10 'Rudolphi, J., Jönsson, M. T., & Gustafsson, L., 2014. ' $false
Append-Run $p10 'Biological legacies buffer local species extinction after logging. ' $true
Append-Run $p10 'Journal of Applied Ecology. 51, 53-62.' $false

# --- fill text for paragraph 11 ---
Append-Run $p11 'Skogsstyrelsen, 2022. ' $false
Append-Run $p11 'Vägledning för hänsyn till knärot. ' $true
Append-Run $p11 'https://www.skogsstyrelsen.se/lag-och-tillsyn/artskydd/vagledningar-och-kunskapsstod-artskydd/vagledning-for-hansyn-till-knarot/' $false

# --- fill text for paragraph 12 ---
Append-Run $p12 'SLU Artdatabanken, 2021. ' $false
Append-Run $p12 'Artfaktablad. Naturvård – artfakta. ' $true
Append-Run $p12 'SLU Artdatabanken, Uppsala ' $false
